$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Player'
$ws.Range("B1").Value = 'Low Price'
$ws.Range("C1").Value = 'High Price'
$ws.Range("A2").Value = 'Reis'
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = '37,600'
$ws.Range("B2").NumberFormat = "General"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = '41,500'
$ws.Range("C2").NumberFormat = "General"
$ws.Range("A3").Value = 'Camacho'
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = '34,200'
$ws.Range("B3").NumberFormat = "General"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = '37,800'
$ws.Range("C3").NumberFormat = "General"
$ws.Range("A4").Value = 'Emil Hansson'
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = '34,000'
$ws.Range("B4").NumberFormat = "General"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = '37,500'
$ws.Range("C4").NumberFormat = "General"
$ws.Range("A5").Value = 'Jony'
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = '35,600'
$ws.Range("B5").NumberFormat = "General"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = '39,300'
$ws.Range("C5").NumberFormat = "General"
$ws.Range("A6").Value = 'Kylian Hazard'
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = '39,500'
$ws.Range("B6").NumberFormat = "General"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = '43,600'
$ws.Range("C6").NumberFormat = "General"
$ws.Range("A7").Value = 'Doğukan Sinik'
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = '38,600'
$ws.Range("B7").NumberFormat = "General"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = '42,700'
$ws.Range("C7").NumberFormat = "General"
$ws.Range("A8").Value = 'Michael Johnston'
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = '39,500'
$ws.Range("B8").NumberFormat = "General"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = '43,600'
$ws.Range("C8").NumberFormat = "General"
$ws.Range("A9").Value = 'Léo Jabá'
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = '38,800'
$ws.Range("B9").NumberFormat = "General"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = '42,900'
$ws.Range("C9").NumberFormat = "General"
$ws.Range("A10").Value = 'Mounir Chouiar'
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = '0'
$ws.Range("B10").NumberFormat = "General"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = '0'
$ws.Range("C10").NumberFormat = "General"
$ws.Range("A11").Value = 'Johan Caballero'
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = '39,500'
$ws.Range("B11").NumberFormat = "General"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = '43,600'
$ws.Range("C11").NumberFormat = "General"
$ws.Range("A12").Value = 'Tobias Mohr'
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = '39,500'
$ws.Range("B12").NumberFormat = "General"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = '43,600'
$ws.Range("C12").NumberFormat = "General"
$ws.Range("A13").Value = 'Song Min Kyu'
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = '39,500'
$ws.Range("B13").NumberFormat = "General"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = '43,600'
$ws.Range("C13").NumberFormat = "General"
$ws.Range("A14").Value = 'Stipe Biuk'
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = '39,500'
$ws.Range("B14").NumberFormat = "General"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = '43,600'
$ws.Range("C14").NumberFormat = "General"
$ws.Range("A15").Value = 'Washington Corozo'
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = '39,500'
$ws.Range("B15").NumberFormat = "General"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = '43,600'
$ws.Range("C15").NumberFormat = "General"
$ws.Range("A16").Value = 'Gustav Mendonca Wikheim'
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = '38,200'
$ws.Range("B16").NumberFormat = "General"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = '42,300'
$ws.Range("C16").NumberFormat = "General"
$ws.Range("A17").Value = 'Octavian Popescu'
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = '39,500'
$ws.Range("B17").NumberFormat = "General"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = '43,600'
$ws.Range("C17").NumberFormat = "General"
$ws.Range("A18").Value = 'Rodrigo Martins'
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = '39,500'
$ws.Range("B18").NumberFormat = "General"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = '43,600'
$ws.Range("C18").NumberFormat = "General"
$ws.Range("A19").Value = 'Lameck Banda'
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = '39,500'
$ws.Range("B19").NumberFormat = "General"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = '43,600'
$ws.Range("C19").NumberFormat = "General"
$ws.Range("A20").Value = 'Carlos Forbs'
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = '39,500'
$ws.Range("B20").NumberFormat = "General"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = '43,600'
$ws.Range("C20").NumberFormat = "General"
$ws.Range("A21").Value = 'Renaldo Cephas'
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = '39,500'
$ws.Range("B21").NumberFormat = "General"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = '43,600'
$ws.Range("C21").NumberFormat = "General"
$ws.Range("A22").Value = 'Joaquín Valiente'
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = '39,500'
$ws.Range("B22").NumberFormat = "General"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = '43,600'
$ws.Range("C22").NumberFormat = "General"
$ws.Range("A23").Value = 'Michael Johnston'
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = '0'
$ws.Range("B23").NumberFormat = "General"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = '0'
$ws.Range("C23").NumberFormat = "General"
